$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 47 and 48 had their match data (columns F..V) swapped back to the
#    original order (Renaissance Zemamra vs Olympique de Safi on row 47,
#    Berkane vs Youssoufia Berrechid on row 48). Columns A..E (Indice, pais,
#    torneio, temporada, data_partida) stay untouched.
# ---------------------------------------------------------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row47 = @{}
$row48 = @{}
foreach ($col in $cols) {
    $row47[$col] = $ws.Range($col + "47").Value()
    $row48[$col] = $ws.Range($col + "48").Value()
}
foreach ($col in $cols) {
    $ws.Range($col + "47").Value = $row48[$col]
    $ws.Range($col + "48").Value = $row47[$col]
}

# ---------------------------------------------------------------------------
# 2) Append two new match rows (63, 64) at the bottom of the sheet, copying
#    the formatting of the last existing row (62) so number formats / bold
#    index styling stay consistent with the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A62:V62").Copy()
$ws.Range("A63:V63").PasteSpecial(-4122)
$ws.Range("A62:V62").Copy()
$ws.Range("A64:V64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# NOTE: this runtime's PowerShell engine does not bind named (-Param value)
# arguments on custom functions, so Set-MatchRow uses positional parameters.
function Set-MatchRow(
    $RowNum, $Indice, $DataPartida, $Home, $HomeGols, $Away, $AwayGols,
    $HomeOpenOdds, $HomeOpenDh, $HomeCloseOdds, $HomeCloseDh,
    $DrawOpenOdds, $DrawOpenDh, $DrawCloseOdds, $DrawCloseDh,
    $AwayOpenOdds, $AwayOpenDh, $AwayCloseOdds, $AwayCloseDh, $Url
) {
    $ws.Range("A$RowNum").Value = $Indice
    $ws.Range("B$RowNum").Value = "morocco"
    $ws.Range("C$RowNum").Value = "botola-pro"
    $ws.Range("D$RowNum").Value = "2023-2024"
    $ws.Range("E$RowNum").Value = $DataPartida
    $ws.Range("F$RowNum").Value = $Home
    $ws.Range("G$RowNum").Value = $HomeGols
    $ws.Range("H$RowNum").Value = $Away
    $ws.Range("I$RowNum").Value = $AwayGols
    $ws.Range("J$RowNum").Value = $HomeOpenOdds
    $ws.Range("K$RowNum").Value = $HomeOpenDh
    $ws.Range("L$RowNum").Value = $HomeCloseOdds
    $ws.Range("M$RowNum").Value = $HomeCloseDh
    $ws.Range("N$RowNum").Value = $DrawOpenOdds
    $ws.Range("O$RowNum").Value = $DrawOpenDh
    $ws.Range("P$RowNum").Value = $DrawCloseOdds
    $ws.Range("Q$RowNum").Value = $DrawCloseDh
    $ws.Range("R$RowNum").Value = $AwayOpenOdds
    $ws.Range("S$RowNum").Value = $AwayOpenDh
    $ws.Range("T$RowNum").Value = $AwayCloseOdds
    $ws.Range("U$RowNum").Value = $AwayCloseDh
    $ws.Range("V$RowNum").Value = $Url
}

Set-MatchRow 63 62 45240.66666666666 "Maghreb Fez" 0 "Union Touarga" 1 `
    1.98 "09/11/2023 04:12" 1.82 "10/11/2023 15:59" `
    2.95 "09/11/2023 04:12" 3.14 "10/11/2023 15:59" `
    3.68 "09/11/2023 04:12" 4.93 "10/11/2023 15:59" `
    "https://www.betexplorer.com/football/morocco/botola-pro/maghreb-fez-union-touarga/tOkxWj89/"

Set-MatchRow 64 63 45240.76041666666 "Moghreb Tetouan" 1 "Mouloudia Oujda" 1 `
    1.86 "09/11/2023 06:42" 1.72 "10/11/2023 18:11" `
    3.09 "09/11/2023 06:42" 3.14 "10/11/2023 18:11" `
    3.93 "09/11/2023 06:42" 5.78 "10/11/2023 18:11" `
    "https://www.betexplorer.com/football/morocco/botola-pro/moghreb-tetouan-mouloudia-oujda/fRgtVANF/"
